# Fruta / hortaliza, semanal
# Inserts two new daily price-report rows (Especial / Segunda) for
# Femacal de La Calera - Frutilla, pushing the existing data block
# (old rows 189-267) down by two rows to (191-269).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right above the current row 189 so the whole
# historical block below shifts down intact (this is what reproduces
# old row N re-appearing at row N+2 with all of its data untouched).
$ws.Rows(189).Insert()
$ws.Rows(189).Insert()

# --- New row 189: "Especial" ---
$ws.Cells.Item(189, 1).Value = 3
$ws.Cells.Item(189, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(189, 3).Value = "Coquimbo"
$ws.Cells.Item(189, 4).Value = 44609
$ws.Cells.Item(189, 5).Value = 5
$ws.Cells.Item(189, 6).Value = "Fruta"
$ws.Cells.Item(189, 7).Value = 100101
$ws.Cells.Item(189, 8).Value = "Berries"
$ws.Cells.Item(189, 9).Value = 100112025
$ws.Cells.Item(189, 10).Value = "Frutilla"
$ws.Cells.Item(189, 11).Value = "Sin especificar"
$ws.Cells.Item(189, 12).Value = "Especial"
$ws.Cells.Item(189, 13).Value = 70
$ws.Cells.Item(189, 14).Value = 6000
$ws.Cells.Item(189, 15).Value = 6000
$ws.Cells.Item(189, 16).Value = 6000
$ws.Cells.Item(189, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(189, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(189, 19).Value = 857
$ws.Cells.Item(189, 20).Value = 7

# --- New row 190: "Segunda" ---
$ws.Cells.Item(190, 1).Value = 3
$ws.Cells.Item(190, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(190, 3).Value = "Coquimbo"
$ws.Cells.Item(190, 4).Value = 44609
$ws.Cells.Item(190, 5).Value = 5
$ws.Cells.Item(190, 6).Value = "Fruta"
$ws.Cells.Item(190, 7).Value = 100101
$ws.Cells.Item(190, 8).Value = "Berries"
$ws.Cells.Item(190, 9).Value = 100112025
$ws.Cells.Item(190, 10).Value = "Frutilla"
$ws.Cells.Item(190, 11).Value = "Sin especificar"
$ws.Cells.Item(190, 12).Value = "Segunda"
$ws.Cells.Item(190, 13).Value = 46
$ws.Cells.Item(190, 14).Value = 4000
$ws.Cells.Item(190, 15).Value = 4000
$ws.Cells.Item(190, 16).Value = 4000
$ws.Cells.Item(190, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(190, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(190, 19).Value = 571
$ws.Cells.Item(190, 20).Value = 7
